$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Escopo")

$ws.Range("F5").Value = 0.8
$ws.Range("E7").Value = "Concluido"
$ws.Range("F7").Value = 1
$ws.Range("F10").Value = 0.6
$ws.Range("F11").Value = 0.6

$ws.Activate()
$ws.Range("F14").Select()
